# Fill in the Snatch / Clean and Jerk weights for the SFM (Semi-Final Men)
# weightlifting sheet, then leave the workbook positioned the way the
# author left it when they saved: SFF's selection parked on F17, and SFM
# as the active sheet/tab with D3 selected.

$wb = $excel.ActiveWorkbook

$wsSFM = $wb.Worksheets.Item("SFM")

# Snatch (C) / Clean and Jerk (D) results for each team, rows 2-7.
$wsSFM.Range("C2").Value = 80
$wsSFM.Range("D2").Value = 120

$wsSFM.Range("C3").Value = 110
$wsSFM.Range("D3").Value = 90

$wsSFM.Range("C4").Value = 105
$wsSFM.Range("D4").Value = 100

$wsSFM.Range("C5").Value = 77.5
$wsSFM.Range("D5").Value = 125

$wsSFM.Range("C6").Value = 92.5
$wsSFM.Range("D6").Value = 110

$wsSFM.Range("C7").Value = 105
$wsSFM.Range("D7").Value = 130

# Leave SFF's own selection parked on F17 before moving away from it.
$wsSFF = $wb.Worksheets.Item("SFF")
$wsSFF.Range("F17").Select()

# Finish on the SFM tab with D3 selected/active, matching the saved file.
$wsSFM.Activate()
$wsSFM.Range("D3").Select()
